# Turn the Chinese header labels on Sheet1 into English ones so the sheet
# can be consumed as a Bootstrap-style table, then spin off a second sheet
# ("Sheet2") that holds the same header plus the first four records as a
# ready-to-use table extract.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- 1. Translate the header row on Sheet1 -------------------------------
$ws1.Range("A1").Value = "name"
$ws1.Range("B1").Value = "department"
$ws1.Range("C1").Value = "speciality"
$ws1.Range("D1").Value = "job"

# --- 2. Add a new worksheet right after Sheet1 ----------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- 3. Populate Sheet2: header + first four data rows from Sheet1 -------
for ($row = 1; $row -le 5; $row++) {
    for ($col = 1; $col -le 4; $col++) {
        $srcCell = $ws1.Cells.Item($row, $col)
        $dstCell = $ws2.Cells.Item($row, $col)
        $dstCell.Value = $srcCell.Value()
    }
}

# --- 4. Restore view/selection state --------------------------------------
# Sheet1: selection spans A1:D5, no longer the active/tab-selected sheet.
$ws1.Range("A1:D5").Select()

# Sheet2: becomes the active sheet with a single-cell selection at F4.
$ws2.Activate()
$ws2.Range("F4").Select()
